# Apply the level1 maze edit:
#  - A handful of special tile markers (P, R, #, S(D), D(L), KW) are no
#    longer used anywhere on the sheet; those cells become plain wall
#    gaps ("_"). This leaves only "W", "_" and "K" in use, so Excel's
#    shared-string table collapses down to just those three values.
#  - The active selection moves from D5 to N7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose special marker is being retired in favor of a plain "_"
$cellsToClear = @("B2", "H2", "J4", "F7", "F9", "C10", "C11", "I11")
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value = "_"
}

# Move the selection/active cell to N7 (was D5)
$ws.Range("N7").Select()
